# Add season record (Wins / Losses / Ties) columns to the PIT_2013 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, reusing the same header style as the existing header row
# (copy the style from the last existing header cell, AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record for every player row (rows 2-51) with the
# team's 2013 season record: 94 wins, 68 losses, 0 ties.
$ws.Range("AD2:AD51").Value = 94
$ws.Range("AE2:AE51").Value = 68
$ws.Range("AF2:AF51").Value = 0
